# Add SQL Connection Logic Dated 09 June 2024
# - Fix the misspelled "Shoping" sheet name to "Shopping"
# - Update the active sheet/selection state left over from the author's
#   last editing session in each sheet

$wb = $excel.ActiveWorkbook

# Correct the sheet name typo: "Shoping" -> "Shopping"
$wsShopping = $wb.Worksheets.Item("Shoping")
$wsShopping.Name = "Shopping"

# The author was last working on the Shopping sheet, with the cursor
# parked at D29 (previously G15) - update the selection there first.
$wsShopping.Activate()
$wsShopping.Range("D29").Select()

# Finish on the Registeration sheet, which becomes the active tab with
# the cursor at I20 (previously the selection was H2:K6 and Login was
# the active tab).
$wsReg = $wb.Worksheets.Item("Registeration")
$wsReg.Activate()
$wsReg.Range("I20").Select()
